# Commit: Sun, May 24, 2020  6:04:56 PM
#
# Re-style the three tables (slides 14-16) that still used the old
# custom "Table_0" style so they use PowerPoint's standard
# "Medium Style 2 - Accent 1" table style instead.

$p = $ppt.ActivePresentation

$oldStyleId = "{9B69767A-36BD-41BD-AE58-8843FAF0CD10}"
$newStyleId = "{C3A4D8BA-4FC4-49BE-AD23-F80122E59843}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
